$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D as text first so numeric-looking price strings
# (e.g. "1.0000", "52.50") are not auto-converted to numbers when assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '23.516.75'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').Value = '1.637.85'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '307.54'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Value = '0.3768'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').Value = '52.50'
$ws.Range('E8').Value = '  +1.78%  '
$ws.Range('D9').Value = '0.3647'
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('D10').Value = '1.268'
$ws.Range('D11').Value = '0.08174'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').Value = '22.97'
$ws.Range('D14').Value = '6.634'
$ws.Range('D15').Value = '0.00001278'
$ws.Range('E15').Value = '  +2.52%  '
$ws.Range('D16').Value = '7.390'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').Value = '1.637.12'
$ws.Range('E17').Value = '  +2.17%  '
$ws.Range('D18').Value = '94.68'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').Value = '0.06935'
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('D20').Value = '18.24'
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('D21').Value = '6.555'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').Value = '1.0000'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '23.531.11'
$ws.Range('E23').Value = '  +1.48%  '
$ws.Range('D24').Value = '12.81'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('D25').Value = '3.090'
$ws.Range('E25').Value = '  +3.49%  '
$ws.Range('D26').Value = '2.425'
$ws.Range('E26').Value = '  +1.31%  '
$ws.Range('D27').Value = '21.26'
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('D28').Value = '151.38'
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').Value = '5.341'
$ws.Range('E29').Value = '  +2.12%  '
$ws.Range('D30').Value = '135.34'
$ws.Range('E30').Value = '  +1.30%  '
$ws.Range('D31').Value = '2.380'
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('D32').Value = '1.821.26'
$ws.Range('E32').Value = '  +2.40%  '
$ws.Range('D33').Value = '6.794'
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('D34').Value = '0.9712'
$ws.Range('E34').Value = '  -0.57%  '
$ws.Range('D35').Value = '0.02821'
$ws.Range('E35').Value = '  +3.70%  '
$ws.Range('D36').Value = '10.30'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').Value = '0.07357'
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('D38').Value = '0.2542'
$ws.Range('E38').Value = '  +1.65%  '
$ws.Range('D39').Value = '6.172'
$ws.Range('E39').Value = '  +0.75%  '
$ws.Range('D40').Value = '0.08874'
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('D41').Value = '1.381'
$ws.Range('E41').Value = '  +1.58%  '
$ws.Range('D42').Value = '0.7097'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').Value = '12.50'
$ws.Range('E43').Value = '  +0.74%  '
$ws.Range('D44').Value = '16.25'
$ws.Range('E44').Value = '  +5.78%  '
$ws.Range('D45').Value = '0.6544'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('D46').Value = '2.337'
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').Value = '4.041'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('D49').Value = '0.07981'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('D50').Value = '129.51'
$ws.Range('E50').Value = '  -2.03%  '
$ws.Range('D51').Value = '1.206'
$ws.Range('E51').Value = '  +0.38%  '

# Restore the original (default) cell style on column D now that the
# text values are safely stored, so no residual style/format attribute
# is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"
